# Refresh the cryptocurrency snapshot: updates each row's "Price" (column D)
# and "Volume(1h)" (column E) figures, and swaps the PaxDollar / mCoin rows
# (41-42), matching the "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $Text into $Cell as literal text. Several of the new "Price" values
# (e.g. "214.65") are syntactically valid numbers; assigning them to .Value
# directly would make Excel silently store them as floating point numbers
# (losing the original text formatting / introducing rounding noise such as
# 214.65000000000001). Prefixing with an apostrophe forces text, exactly as
# typing e.g. `'214.65` into Excel would.
function Set-TextValue {
    param($Cell, [string]$Text)
    if ($Text -match '^\s*[-+]?\d+(\.\d+)?\s*$') {
        $ws.Range($Cell).Value = "'" + $Text
    } else {
        $ws.Range($Cell).Value = $Text
    }
}

# Row 2: D2='25.841.69', E2='  -1.20%  '
Set-TextValue "D2" '25.841.69'
Set-TextValue "E2" '  -1.20%  '

# Row 3: D3='1.633.15', E3='  -1.25%  '
Set-TextValue "D3" '1.633.15'
Set-TextValue "E3" '  -1.25%  '

# Row 4: E4='  -0.36%  '
Set-TextValue "E4" '  -0.36%  '

# Row 5: D5='214.65', E5='  -0.39%  '
Set-TextValue "D5" '214.65'
Set-TextValue "E5" '  -0.39%  '

# Row 6: D6='0.5018', E6='  -1.83%  '
Set-TextValue "D6" '0.5018'
Set-TextValue "E6" '  -1.83%  '

# Row 7: E7='  -0.41%  '
Set-TextValue "E7" '  -0.41%  '

# Row 8: D8='0.2566', E8='  -0.81%  '
Set-TextValue "D8" '0.2566'
Set-TextValue "E8" '  -0.81%  '

# Row 9: D9='0.06388', E9='  -0.24%  '
Set-TextValue "D9" '0.06388'
Set-TextValue "E9" '  -0.24%  '

# Row 10: E10='  -1.84%  '
Set-TextValue "E10" '  -1.84%  '

# Row 11: D11='0.07685', E11='  -1.65%  '
Set-TextValue "D11" '0.07685'
Set-TextValue "E11" '  -1.65%  '

# Row 12: D12='1.636.60', E12='  -1.10%  '
Set-TextValue "D12" '1.636.60'
Set-TextValue "E12" '  -1.10%  '

# Row 13: D13='4.240', E13='  -1.04%  '
Set-TextValue "D13" '4.240'
Set-TextValue "E13" '  -1.04%  '

# Row 14: D14='1.858.47', E14='  -1.27%  '
Set-TextValue "D14" '1.858.47'
Set-TextValue "E14" '  -1.27%  '

# Row 15: D15='0.5427', E15='  -1.88%  '
Set-TextValue "D15" '0.5427'
Set-TextValue "E15" '  -1.88%  '

# Row 16: D16='0.0₅7915', E16='  -1.31%  '
Set-TextValue "D16" '0.0₅7915'
Set-TextValue "E16" '  -1.31%  '

# Row 17: D17='63.38', E17='  -0.83%  '
Set-TextValue "D17" '63.38'
Set-TextValue "E17" '  -0.83%  '

# Row 18: D18='25.845.93', E18='  -1.27%  '
Set-TextValue "D18" '25.845.93'
Set-TextValue "E18" '  -1.27%  '

# Row 19: E19='  -0.34%  '
Set-TextValue "E19" '  -0.34%  '

# Row 20: D20='201.65', E20='  -3.51%  '
Set-TextValue "D20" '201.65'
Set-TextValue "E20" '  -3.51%  '

# Row 21: D21='4.321', E21='  -2.09%  '
Set-TextValue "D21" '4.321'
Set-TextValue "E21" '  -2.09%  '

# Row 22: D22='9.924', E22='  -1.42%  '
Set-TextValue "D22" '9.924'
Set-TextValue "E22" '  -1.42%  '

# Row 23: D23='5.969', E23='  -1.00%  '
Set-TextValue "D23" '5.969'
Set-TextValue "E23" '  -1.00%  '

# Row 24: E24='  -0.28%  '
Set-TextValue "E24" '  -0.28%  '

# Row 25: E25='  +10.43%  '
Set-TextValue "E25" '  +10.43%  '

# Row 26: D26='141.02', E26='  -1.43%  '
Set-TextValue "D26" '141.02'
Set-TextValue "E26" '  -1.43%  '

# Row 27: E27='  -2.45%  '
Set-TextValue "E27" '  -2.45%  '

# Row 28: E28='  -0.69%  '
Set-TextValue "E28" '  -0.69%  '

# Row 29: D29='6.692', E29='  -3.99%  '
Set-TextValue "D29" '6.692'
Set-TextValue "E29" '  -3.99%  '

# Row 30: E30='  -0.44%  '
Set-TextValue "E30" '  -0.44%  '

# Row 31: E31='  -2.98%  '
Set-TextValue "E31" '  -2.98%  '

# Row 32: E32='  -2.93%  '
Set-TextValue "E32" '  -2.93%  '

# Row 33: D33='3.174', E33='  -1.49%  '
Set-TextValue "D33" '3.174'
Set-TextValue "E33" '  -1.49%  '

# Row 34: D34='1.536', E34='  -1.58%  '
Set-TextValue "D34" '1.536'
Set-TextValue "E34" '  -1.58%  '

# Row 35: D35='2.365', E35='  -0.38%  '
Set-TextValue "D35" '2.365'
Set-TextValue "E35" '  -0.38%  '

# Row 36: D36='1.169.03', E36='  +1.13%  '
Set-TextValue "D36" '1.169.03'
Set-TextValue "E36" '  +1.13%  '

# Row 37: D37='0.8908', E37='  -4.14%  '
Set-TextValue "D37" '0.8908'
Set-TextValue "E37" '  -4.14%  '

# Row 38: D38='2.619', E38='  -4.46%  '
Set-TextValue "D38" '2.619'
Set-TextValue "E38" '  -4.46%  '

# Row 39: D39='0.5596', E39='  -1.43%  '
Set-TextValue "D39" '0.5596'
Set-TextValue "E39" '  -1.43%  '

# Row 40: E40='  -1.94%  '
Set-TextValue "E40" '  -1.94%  '

# Row 41: B41='mCoin', C41='https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin', D41='2.552', E41='  -0.26%  '
Set-TextValue "B41" 'mCoin'
Set-TextValue "C41" 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
Set-TextValue "D41" '2.552'
Set-TextValue "E41" '  -0.26%  '

# Row 42: B42='PaxDollar', C42='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', D42='1.002', E42='  -0.38%  '
Set-TextValue "B42" 'PaxDollar'
Set-TextValue "C42" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D42" '1.002'
Set-TextValue "E42" '  -0.38%  '

# Row 43: D43='5.680', E43='  +0.61%  '
Set-TextValue "D43" '5.680'
Set-TextValue "E43" '  +0.61%  '

# Row 44: D44='0.8068', E44='  -3.21%  '
Set-TextValue "D44" '0.8068'
Set-TextValue "E44" '  -3.21%  '

# Row 45: D45='99.33', E45='  -0.57%  '
Set-TextValue "D45" '99.33'
Set-TextValue "E45" '  -0.57%  '

# Row 46: D46='1.770.76', E46='  -1.27%  '
Set-TextValue "D46" '1.770.76'
Set-TextValue "E46" '  -1.27%  '

# Row 47: E47='  -1.25%  '
Set-TextValue "E47" '  -1.25%  '

# Row 48: E48='  -0.72%  '
Set-TextValue "E48" '  -0.72%  '

# Row 49: D49='1.003', E49='  +0.08%  '
Set-TextValue "D49" '1.003'
Set-TextValue "E49" '  +0.08%  '

# Row 50: D50='54.69', E50='  -1.69%  '
Set-TextValue "D50" '54.69'
Set-TextValue "E50" '  -1.69%  '

# Row 51: D51='0.05068', E51='  +0.52%  '
Set-TextValue "D51" '0.05068'
Set-TextValue "E51" '  +0.52%  '
